# Zeiterfassung: SOLL eingetragen -> IST jeweils am Abend eintragen
#
# For a number of task rows on the "Zeitplanung" sheet, the value that had
# been entered into one of the daily-tracking cells (columns G:BJ, the "IST"
# grid) is moved into column C (the "SOLL" cell for that task), and the
# daily-tracking cell that used to hold it is cleared out. Column D and the
# various SUM() totals are formulas and recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitplanung")

# --- Set new SOLL (column C) values -----------------------------------
$ws.Range("C12").Value = 2.5
$ws.Range("C19").Value = 5
$ws.Range("C20").Value = 4
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 4
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 9
$ws.Range("C25").Value = 4
$ws.Range("C26").Value = 4
$ws.Range("C27").Value = 2
$ws.Range("C28").Value = 2
$ws.Range("C29").Value = 3
$ws.Range("C35").Value = 1
$ws.Range("C36").Value = 2
$ws.Range("C37").Value = 3
$ws.Range("C42").Value = 1
$ws.Range("C43").Value = 1

# --- Clear the daily IST-grid cells whose values were moved above ------
$ws.Range("K12").ClearContents()
$ws.Range("R12").ClearContents()
$ws.Range("Y12").ClearContents()
$ws.Range("AD12").ClearContents()
$ws.Range("BA12").ClearContents()
$ws.Range("W21").ClearContents()
$ws.Range("W22").ClearContents()
$ws.Range("X23").ClearContents()
$ws.Range("Y23").ClearContents()
$ws.Range("AB24").ClearContents()
$ws.Range("AC24").ClearContents()
$ws.Range("Y25").ClearContents()
$ws.Range("Y26").ClearContents()
$ws.Range("AB26").ClearContents()
$ws.Range("AC27").ClearContents()
$ws.Range("AC28").ClearContents()
$ws.Range("AC29").ClearContents()
$ws.Range("AD29").ClearContents()
$ws.Range("AD35").ClearContents()
$ws.Range("AD36").ClearContents()
$ws.Range("AD37").ClearContents()
$ws.Range("BA39").ClearContents()
$ws.Range("BA42").ClearContents()
$ws.Range("BA43").ClearContents()

$excel.ActiveWorkbook.Save()
